$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting all existing data right by one.
$ws.Range("A:A").EntireColumn.Insert()

# Header for the new column.
$ws.Range("A1").Value() = "Match ID"

# Bold the new column for the visible rows (1-19); this creates the
# font-only style (bold, no border/alignment) used for the Match ID column.
$ws.Range("A1:A19").Font.Bold() = $true

# Fill in the Match ID value for every data row (rows 4-20; row 20 is the
# hidden summary/total row).
for ($r = 4; $r -le 20; $r++) {
  $ws.Cells.Item($r, 1).Value() = 17
}

# Restore the sheet selection.
$ws.Range("A1:A19").Select()
